$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.612.26"
$ws.Range("E2").Value = "  +3.80%  "

$ws.Range("D3").Value = "3.505.63"
$ws.Range("E3").Value = "  +2.38%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'591.02"
$ws.Range("E5").Value = "  +3.35%  "

$ws.Range("D6").Value = "'170.05"
$ws.Range("E6").Value = "  +5.38%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "3.503.22"
$ws.Range("E8").Value = "  +2.28%  "

$ws.Range("D9").Value = "'0.593"
$ws.Range("E9").Value = "  +6.87%  "

$ws.Range("D10").Value = "'7.35"
$ws.Range("E10").Value = "  +0.85%  "

$ws.Range("E11").Value = "  +4.84%  "

$ws.Range("D12").Value = "'0.440"
$ws.Range("E12").Value = "  +3.72%  "

$ws.Range("D13").Value = "4.113.21"
$ws.Range("E13").Value = "  +2.39%  "

$ws.Range("E14").Value = "  -0.49%  "

$ws.Range("D15").Value = "'28.43"
$ws.Range("E15").Value = "  +4.93%  "

$ws.Range("E16").Value = "  +2.29%  "

$ws.Range("D17").Value = "66.659.22"
$ws.Range("E17").Value = "  +3.82%  "

$ws.Range("D18").Value = "3.511.92"
$ws.Range("E18").Value = "  +2.20%  "

$ws.Range("E19").Value = "  +4.39%  "

$ws.Range("D20").Value = "'14.13"
$ws.Range("E20").Value = "  +3.99%  "

$ws.Range("D21").Value = "'390.25"
$ws.Range("E21").Value = "  +3.20%  "

$ws.Range("E22").Value = "  +2.04%  "

$ws.Range("D23").Value = "'73.09"
$ws.Range("E23").Value = "  +2.25%  "

$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("D25").Value = "'0.538"
$ws.Range("E25").Value = "  +3.98%  "

$ws.Range("E26").Value = "  +5.61%  "

$ws.Range("D27").Value = "'10.36"
$ws.Range("E27").Value = "  +8.53%  "

$ws.Range("E28").Value = "  +2.47%  "

$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "'1.50"
$ws.Range("E30").Value = "  +7.17%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'6.37"
$ws.Range("E31").Value = "  +5.82%  "

$ws.Range("D33").Value = "'23.64"
$ws.Range("E33").Value = "  +2.98%  "

$ws.Range("D34").Value = "'7.42"
$ws.Range("E34").Value = "  +5.01%  "

$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("E36").Value = "  +7.32%  "

$ws.Range("D37").Value = "'162.68"
$ws.Range("E37").Value = "  +1.92%  "

$ws.Range("D38").Value = "'0.885"
$ws.Range("E38").Value = "  +3.04%  "

$ws.Range("E39").Value = "  +4.66%  "

$ws.Range("E41").Value = "  +5.14%  "

$ws.Range("E42").Value = "  +2.88%  "

$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'27.54"
$ws.Range("E43").Value = "  +5.41%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'26.56"
$ws.Range("E44").Value = "  +2.96%  "

$ws.Range("D45").Value = "2.812.49"
$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("D46").Value = "'43.20"
$ws.Range("E46").Value = "  +0.45%  "

$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "'2.53"
$ws.Range("E47").Value = "  +4.96%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0312"
$ws.Range("E48").Value = "  +2.25%  "

$ws.Range("D49").Value = "'355.05"
$ws.Range("E49").Value = "  +5.02%  "

$ws.Range("E50").Value = "  +3.43%  "

$ws.Range("D51").Value = "'33.81"
$ws.Range("E51").Value = "  +12.90%  "
